# Update "想去人数" (F column) values on sheets "展览" and "全部类型"
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    "F2" = 734
    "F3" = 480
    "F6" = 618
    "F7" = 38
    "F8" = 335
    "F10" = 396
    "F12" = 116
    "F17" = 1060
    "F20" = 351
    "F21" = 184
    "F23" = 158
    "F25" = 92
    "F26" = 242
    "F27" = 271
    "F29" = 1655
    "F35" = 3857
    "F37" = 447
    "F38" = 222
    "F39" = 971
    "F40" = 79
    "F43" = 86
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}

Write-Output "Updated F-column values on sheets: $($sheetNames -join ', ')"
